$d = $word.ActiveDocument

# Merge "H. R. " + "Allcock" + ", " (split across runs with proofErr
# spell-check wrappers) into a single run's text "H. R. Allcock, ".
$d.Content.Find.Execute(
    "H. R. Allcock, ", $true, $false, $false, $false, $false,
    $true, 1, $false, "H. R. Allcock, ", 2
)

# Merge ", Vol. 2 (Eds: H. " + "Baltes" + ", W. " (also split across
# runs with proofErr spell-check wrappers) into a single run.
$d.Content.Find.Execute(
    ", Vol. 2 (Eds: H. Baltes, W. ", $true, $false, $false, $false, $false,
    $true, 1, $false, ", Vol. 2 (Eds: H. Baltes, W. ", 2
)
